$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Country: Afghanistan -> India
$ws.Range("D2").Value = "India"
$ws.Range("D3").Value = "India"

# State: Goa -> Gujarat
$ws.Range("E2").Value = "Gujarat"
$ws.Range("E3").Value = "Gujarat"

# New Telephone number for row 2
$ws.Range("H2").Value = 7418521478

# Re-fit the State/City columns so the new, wider values ("Gujarat") are fully visible
$ws.Columns("E:F").AutoFit()
$ws.Columns("E").ColumnWidth = 6.666666666666667
$ws.Columns("F").ColumnWidth = 10.833333333333334

# Move the active selection to H2
$ws.Range("H2").Select()
